$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert three new paragraphs right after the
#    "- wyświetlanie listy wszystkich obiektów 10.11.2016" paragraph
#    (i.e. right before the "-System komentowania/wyszukiwania/oceniania"
#    paragraph):
#       - ogólny wygląd strony
#       <tab>- strona główna  18.11.2016
#       <lastRenderedPageBreak><tab>-menu  18.11.2016
# ---------------------------------------------------------------------
$anchorIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*wyświetlanie listy wszystkich*") {
        $anchorIndex = $i
    }
}

$r = $d.Paragraphs($anchorIndex).Range
$r.Collapse(0)
$r.InsertParagraphAfter()
$d.Paragraphs($anchorIndex + 1).Range.InsertXML(
    "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t>- ogólny wygląd strony</w:t></w:r></w:p>"
)

$r2 = $d.Paragraphs($anchorIndex + 1).Range
$r2.Collapse(0)
$r2.InsertParagraphAfter()
$d.Paragraphs($anchorIndex + 2).Range.InsertXML(
    "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:tab/><w:t>- strona główna  18.11.2016</w:t></w:r></w:p>"
)

$r3 = $d.Paragraphs($anchorIndex + 2).Range
$r3.Collapse(0)
$r3.InsertParagraphAfter()
$d.Paragraphs($anchorIndex + 3).Range.InsertXML(
    "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:lastRenderedPageBreak/><w:tab/><w:t>-menu  18.11.2016</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# 2. Drop the lastRenderedPageBreak that used to sit on the "-ocenianie"
#    paragraph (it effectively moved up onto the new "-menu" paragraph
#    above).
# ---------------------------------------------------------------------
$ocenIndex = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "*-ocenianie*") {
        $ocenIndex = $i
    }
}
$d.Paragraphs($ocenIndex).Range.InsertXML(
    "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:tab/><w:t>-ocenianie</w:t></w:r><w:r><w:tab/><w:t>25.11.2016</w:t></w:r></w:p>"
)

# ---------------------------------------------------------------------
# 3. Remove the obsolete scratch notes at the end of the document:
#       "logowanie przesunać niżej "
#       (empty paragraph)
#       (empty paragraph)
#       "dodawanie edytowanie"
#    while keeping the following empty paragraph and the final
#    "usuwanie z customowym komunikatem..." paragraph untouched.
# ---------------------------------------------------------------------
$startIdx = $null
$endIdx = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs($i).Range.Text
    if ($t -like "*logowanie przesunać niżej*") {
        $startIdx = $i
    }
    if ($t -like "*dodawanie edytowanie*") {
        $endIdx = $i
    }
}
$startRange = $d.Paragraphs($startIdx).Range
$endRange = $d.Paragraphs($endIdx).Range
$delRange = $d.Range($startRange.Start, $endRange.End)
$delRange.Delete()
